$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Fix the age-15up late progression rate formula to match the age0to5 / age5to15 rows
$ws.Range("B35").Formula = "=6.8/1000000*365"

# Insert a new row at 91 (pushing existing row 91 and below down by one)
$ws.Rows("91:91").Insert(1)

# Copy the formatting of the row above (row 90, prison_age_min) into the new row 91
$ws.Rows("90:90").Copy()
$ws.Rows("91:91").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new constant: riskgroup_startage_dorm = 15
$ws.Range("A91").Value = "riskgroup_startage_dorm"
$ws.Range("B91").Value = 15
$ws.Range("C91").Value = $null
$ws.Range("D91").Value = $null
$ws.Range("E91").Value = $null

# Make the constants sheet the active tab / sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 72
$ws.Range("D91").Select()

# time_variants sheet should no longer be the tab shown as selected
$ws2 = $wb.Worksheets.Item("time_variants")
$ws2.Select()
$ws.Select()
